# Apply the edits described by the diff:
#  - column G (7) width shrinks from 3.140625 to 2.140625
#  - column I (9) and J (10) widths grow from 3.140625 to 5.7109375
#  - column K (11) width shrinks from 7.7109375 to 5.7109375
#  - several values in row 1 are updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes -------------------------------------------------
# Excel's COM ColumnWidth property is expressed in "characters" and is
# converted internally to the OOXML <col width="..."/> units; the values
# below are the character-width inputs that land closest to the target
# OOXML widths taken from the diff (2.140625 and 5.7109375).
$ws.Columns.Item(7).ColumnWidth = 1.333333    # -> width ~2.140625 (was 3.140625)
$ws.Columns.Item(9).ColumnWidth = 4.833333    # -> width ~5.7109375 (was 3.140625)
$ws.Columns.Item(10).ColumnWidth = 4.833333   # -> width ~5.7109375 (was 3.140625)
$ws.Columns.Item(11).ColumnWidth = 4.833333   # -> width ~5.7109375 (was 7.7109375)

# --- Cell value changes (row 1) --------------------------------------------
$ws.Range("A1").Value = 3          # was 4
$ws.Range("C1").Value = 10         # was 15
$ws.Range("E1").Value = 30         # was 11
$ws.Range("G1").Value = 6          # was 10
$ws.Range("H1").Value = 31         # was 24
$ws.Range("I1").Value = 0.065      # was 18
$ws.Range("J1").Value = 0.061      # was 30
$ws.Range("K1").Value = 0.025      # was 0.02299
